$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# List of cell updates: row, column index, new text value
# Values are written as text (NumberFormat "@") to preserve the original
# string representation (e.g. "1.001", "29.444.26"), then formatting is
# cleared so no stray style/number-format is left on the cell.
$updates = @(
    ,@(2, 4, '29.444.26')
    ,@(2, 5, '  +0.06%  ')
    ,@(3, 4, '1.853.21')
    ,@(3, 5, '  +0.23%  ')
    ,@(4, 4, '1.001')
    ,@(4, 5, '  +0.11%  ')
    ,@(5, 4, '241.33')
    ,@(5, 5, '  +0.33%  ')
    ,@(6, 4, '0.6324')
    ,@(6, 5, '  +0.83%  ')
    ,@(7, 5, '  +0.24%  ')
    ,@(8, 4, '0.07536')
    ,@(8, 5, '  -2.11%  ')
    ,@(9, 4, '0.2920')
    ,@(9, 5, '  +0.23%  ')
    ,@(10, 4, '24.48')
    ,@(10, 5, '  -1.21%  ')
    ,@(11, 4, '0.07757')
    ,@(11, 5, '  +0.09%  ')
    ,@(12, 4, '1.854.16')
    ,@(12, 5, '  -0.44%  ')
    ,@(13, 4, '5.031')
    ,@(13, 5, '  +0.14%  ')
    ,@(14, 4, '0.6828')
    ,@(14, 5, '  +0.23%  ')
    ,@(15, 4, '0.00001041')
    ,@(15, 5, '  -3.47%  ')
    ,@(16, 4, '83.27')
    ,@(16, 5, '  -0.32%  ')
    ,@(17, 4, '2.120.53')
    ,@(17, 5, '  +0.59%  ')
    ,@(18, 4, '6.135')
    ,@(18, 5, '  -0.53%  ')
    ,@(19, 4, '29.470.76')
    ,@(19, 5, '  +0.09%  ')
    ,@(20, 4, '229.61')
    ,@(20, 5, '  +0.49%  ')
    ,@(21, 4, '12.37')
    ,@(21, 5, '  -0.16%  ')
    ,@(22, 5, '  +0.20%  ')
    ,@(23, 4, '7.474')
    ,@(23, 5, '  +0.76%  ')
    ,@(24, 4, '1.004')
    ,@(24, 5, '  +0.33%  ')
    ,@(25, 4, '159.46')
    ,@(25, 5, '  +1.46%  ')
    ,@(26, 4, '0.1388')
    ,@(26, 5, '  +1.04%  ')
    ,@(27, 4, '8.438')
    ,@(27, 5, '  +0.47%  ')
    ,@(28, 4, '17.65')
    ,@(28, 5, '  -0.31%  ')
    ,@(29, 4, '1.421')
    ,@(29, 5, '  +5.76%  ')
    ,@(30, 4, '1.477')
    ,@(30, 5, '  +0.85%  ')
    ,@(31, 4, '0.05702')
    ,@(31, 5, '  +1.04%  ')
    ,@(32, 4, '4.141')
    ,@(32, 5, '  +0.55%  ')
    ,@(33, 4, '4.052')
    ,@(33, 5, '  +0.46%  ')
    ,@(34, 4, '1.157')
    ,@(34, 5, '  -0.43%  ')
    ,@(35, 4, '1.818')
    ,@(35, 5, '  -1.22%  ')
    ,@(36, 4, '0.6985')
    ,@(36, 5, '  -1.40%  ')
    ,@(37, 4, '2.593')
    ,@(37, 5, '  +0.00%  ')
    ,@(38, 4, '2.842')
    ,@(38, 5, '  +2.67%  ')
    ,@(39, 4, '1.252.78')
    ,@(39, 5, '  +1.99%  ')
    ,@(40, 4, '0.01832')
    ,@(40, 5, '  +2.44%  ')
    ,@(41, 5, '  +0.82%  ')
    ,@(42, 4, '0.9089')
    ,@(42, 5, '  +0.46%  ')
    ,@(43, 4, '1.002')
    ,@(43, 5, '  +0.09%  ')
    ,@(44, 4, '2.019.85')
    ,@(44, 5, '  +0.24%  ')
    ,@(45, 4, '101.70')
    ,@(45, 5, '  +0.02%  ')
    ,@(46, 4, '66.05')
    ,@(46, 5, '  +0.32%  ')
    ,@(47, 4, '7.111')
    ,@(47, 5, '  -0.75%  ')
    ,@(48, 2, 'BabyDogeCoin')
    ,@(48, 3, 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge')
    ,@(48, 4, '0.00000000117')
    ,@(48, 5, '  -1.93%  ')
    ,@(49, 2, 'Algorand')
    ,@(49, 3, 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo')
    ,@(49, 4, '0.1163')
    ,@(49, 5, '  +0.49%  ')
    ,@(50, 4, '9.031')
    ,@(50, 5, '  +0.34%  ')
    ,@(51, 4, '0.3966')
    ,@(51, 5, '  -1.07%  ')
)

foreach ($u in $updates) {
    $r = $u[0]
    $c = $u[1]
    $val = $u[2]
    $cell = $ws.Cells.Item($r, $c)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}